$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new team-name labels first, in the same order the shared
# string table grows in the target file (Muchacho1, Bonjour1, Bonjour2,
# Muchacho2), so the resulting sharedStrings.xml ordering matches.
$ws.Range("A6").Value = "Muchacho1"
$ws.Range("A4").Value = "Bonjour1"
$ws.Range("A5").Value = "Bonjour2"
$ws.Range("A7").Value = "Muchacho2"

# Fill in the numeric score columns for the new rows (4-7)
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1

$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 2

$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3

$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 4

# Update the selected cell to match the author's final cursor position
$ws.Range("G9").Select()
